# LNI-324: Added test coverage of 'subsubject', and 'approval' and
# 'laidInDraft' rubrics.
#
# 1. Insert a new "Subsub" paragraph ("superannuation") right after the
#    "subject" paragraph ("NATIONAL HEALTH SERVICE") and before the
#    "Title" paragraph.
# 2. Insert new "Approval" ("Approved by both Houses of Parliament") and
#    "LaidDraft" ("Laid before Parliament in draft") paragraphs right
#    after the "Title" paragraph and before the "Made" paragraph.
#
# Paragraph references returned by $d.Paragraphs.Item($i) track the
# $i'th paragraph slot rather than a stable node identity, so any
# insertion shifts the index of everything after it. To keep every
# lookup valid we resolve indices by style name right before each use,
# and we perform the edits from the bottom of the document upwards so
# earlier (already-used) indices are never invalidated by later work.

$d = $word.ActiveDocument

function Get-ParaIndexByStyle([string]$styleName) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Style.NameLocal -eq $styleName) {
            return $i
        }
    }
    return -1
}

# --- 2. "Approval" / "LaidDraft" after "Title" (handled first: it is the
#        later insertion point, so doing it before the "Subsub" edit
#        keeps the "Title" paragraph's index stable while we use it). ---
$titleIndex = Get-ParaIndexByStyle "Title"
$titlePara = $d.Paragraphs.Item($titleIndex)
$titlePara.Range.InsertParagraphAfter()

$approvalPara = $d.Paragraphs.Item($titleIndex + 1)
$approvalPara.Style = "Approval"
$approvalPara.Range.Text = "Approved by both Houses of Parliament"

$approvalPara.Range.InsertParagraphAfter()

$laidDraftPara = $d.Paragraphs.Item($titleIndex + 2)
$laidDraftPara.Style = "LaidDraft"
$laidDraftPara.Range.Text = "Laid before Parliament in draft"

# --- 1. "Subsub" after "subject". ---
$subjectIndex = Get-ParaIndexByStyle "subject"
$subjectPara = $d.Paragraphs.Item($subjectIndex)
$subjectPara.Range.InsertParagraphAfter()

$subsubPara = $d.Paragraphs.Item($subjectIndex + 1)
$subsubPara.Style = "Subsub"
$subsubPara.Range.Text = "superannuation"
